$d = $word.ActiveDocument

# Word constants used below:
#   wdFindContinue      = 1
#   wdReplaceOne        = 1
#   FindContinue (Wrap) = 1

# 1) First "Data" line ("Data: 09/10/2024" -> "Data: 11/10/2024").
#    Only the day portion ("09") changes to "11"; the rest of the
#    paragraph (month "/10", year "/2024", and the bold "Data" label)
#    is left untouched.
$p1 = $d.Paragraphs(3).Range
$p1.Find.Execute("09", $false, $false, $false, $false, $false, $true, 1, $false, "11", 1)

# 2) "Próxima Reunião" -> "Data" line ("Data: 10/10/2024" -> "Data: 13/10/2024").
#    Only the first occurrence (day portion "10") changes to "13"; the
#    month "/10" and year "/2024" stay the same.
$p2 = $d.Paragraphs(27).Range
$p2.Find.Execute("10", $false, $false, $false, $false, $false, $true, 1, $false, "13", 1)
